$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.414.83'
$ws.Range("E2").Value = '  +3.56%  '

# Row 3
$ws.Range("D3").Value = '2.308.99'
$ws.Range("E3").Value = '  +2.78%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.62'
$ws.Range("E5").Value = '  +0.85%  '

# Row 6
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.11'
$ws.Range("E6").Value = '  +9.40%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("E7").Value = '  +0.70%  '

# Row 8
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  +5.97%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.84'
$ws.Range("E10").Value = '  +3.66%  '

# Row 11
$ws.Range("E11").Value = '  +2.61%  '

# Row 12
$ws.Range("E12").Value = '  -0.12%  '

# Row 13
$ws.Range("E13").Value = '  -1.03%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.96'
$ws.Range("E14").Value = '  +3.08%  '

# Row 15
$ws.Range("D15").Value = '2.673.45'
$ws.Range("E15").Value = '  +3.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.10'
$ws.Range("E16").Value = '  +5.36%  '

# Row 17
$ws.Range("D17").Value = '2.314.49'
$ws.Range("E17").Value = '  +3.61%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.801'
$ws.Range("E18").Value = '  +2.83%  '

# Row 19
$ws.Range("D19").Value = '43.357.27'
$ws.Range("E19").Value = '  +3.59%  '

# Row 20
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0922'
$ws.Range("E20").Value = '  +2.53%  '

# Row 21
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.83'
$ws.Range("E21").Value = '  -2.43%  '

# Row 22
$ws.Range("E22").Value = '  +4.90%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.93'
$ws.Range("E23").Value = '  +1.25%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '240.19'
$ws.Range("E24").Value = '  +2.19%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("E25").Value = '  +4.02%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.60'
$ws.Range("E26").Value = '  +1.46%  '

# Row 27
$ws.Range("E27").Value = '  +0.12%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.06'
$ws.Range("E28").Value = '  +8.24%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("E29").Value = '  +4.71%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.17'
$ws.Range("E30").Value = '  -4.22%  '

# Row 31
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.58'
$ws.Range("E31").Value = '  +1.46%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '162.35'
$ws.Range("E32").Value = '  -2.89%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.23'
$ws.Range("E33").Value = '  +1.76%  '

# Row 34
$ws.Range("E34").Value = '  -0.08%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.30'
$ws.Range("E35").Value = '  +5.27%  '

# Row 36
$ws.Range("E36").Value = '  +6.40%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0731'
$ws.Range("E37").Value = '  +2.02%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.59'
$ws.Range("E38").Value = '  +14.17%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").Value = '  -2.37%  '

# Row 40
$ws.Range("E40").Value = '  +3.17%  '

# Row 41
$ws.Range("E41").Value = '  +3.95%  '

# Row 42
$ws.Range("E42").Value = '  +0.53%  '

# Row 43
$ws.Range("E43").Value = '  +14.80%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0288'
$ws.Range("E44").Value = '  +3.06%  '

# Row 45
$ws.Range("D45").Value = '1.963.58'
$ws.Range("E45").Value = '  +1.43%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.79'
$ws.Range("E46").Value = '  +2.15%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.05'
$ws.Range("E47").Value = '  +6.12%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.26'
$ws.Range("E48").Value = '  +6.72%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '57.80'
$ws.Range("E49").Value = '  +7.63%  '

# Row 50
$ws.Range("E50").Value = '  +1.54%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.57'
$ws.Range("E51").Value = '  +8.06%  '
